# The workbook's two "full listing" sheets (展览 and 全部类型) both start with
# an event row for "南宁·AP动漫游戏嘉年华" (2024-01-27) that was removed from the
# upstream feed. Every following event shifted up one row and a fresh event
# was appended at the bottom; a handful of numeric "want-to-go" counts were
# also refreshed, and the first remaining event's lowest-price column became
# unavailable ("不可售") instead of a numeric price.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Drop the old top data row (2024-01-27 · 南宁·AP动漫游戏嘉年华); this shifts
    # every subsequent row up by one and removes the old trailing row, so the
    # used range shrinks from A1:I7 to A1:I6 automatically.
    $ws.Rows.Item(2).Delete()

    # Renumber the leading index column back to 1..5.
    $ws.Range("A2").Value = 1
    $ws.Range("A3").Value = 2
    $ws.Range("A4").Value = 3
    $ws.Range("A5").Value = 4
    $ws.Range("A6").Value = 5

    # The event that is now first ("南宁·第一届异次元动漫嘉年华") lost ticket
    # availability - its lowest-price cell becomes a text marker.
    $ws.Range("G2").Value = "不可售"

    # Refreshed "want to go" counts for three of the surviving events.
    $ws.Range("F3").Value = 1674
    $ws.Range("F4").Value = 7751
    $ws.Range("F6").Value = 224
}
